$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.853.27"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "'1.782.35"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'312.01"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5113"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "'0.3760"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "'0.07779"
$ws.Range("E9").Value = "  -7.74%  "
$ws.Range("D10").Value = "'41.23"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "'1.084"
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "'6.183"
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("D14").Value = "'20.15"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").Value = "'1.779.81"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "'7.167"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("D17").Value = "'91.73"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "'0.00001070"
$ws.Range("E18").Value = "  -5.96%  "
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'16.99"
$ws.Range("E21").Value = "  -4.12%  "
$ws.Range("D22").Value = "'5.909"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "'27.915.99"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").Value = "'10.91"
$ws.Range("E24").Value = "  -4.60%  "
$ws.Range("D25").Value = "'2.246"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "'158.24"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").Value = "'1.986.59"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "'2.348"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").Value = "'122.24"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "'0.1078"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  -5.30%  "
$ws.Range("D33").Value = "'3.628"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D35").Value = "'0.07084"
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").Value = "'0.2122"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.001"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.46"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.509"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("D43").Value = "'1.149"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").Value = "'1.327"
$ws.Range("E44").Value = "  -5.42%  "
$ws.Range("D45").Value = "'13.03"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("D46").Value = "'0.5923"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "'3.720"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "'126.23"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'1.211"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("E50").Value = "  -4.86%  "
$ws.Range("D51").Value = "'0.06717"
$ws.Range("E51").Value = "  -3.78%  "
